$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.237.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.53%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.285.28'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.80%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.60%  '

$ws.Range('E6').Value = '  -2.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.621'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.96%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.604'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.67%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.32'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.29%  '

$ws.Range('E11').Value = '  -0.83%  '

$ws.Range('E12').Value = '  -2.31%  '

$ws.Range('E13').Value = '  -1.16%  '

$ws.Range('E14').Value = '  -1.71%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.22%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.632.38'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.280.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.13%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.322.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.76%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.95%  '

$ws.Range('E20').Value = '  -0.80%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +27.30%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.36%  '

$ws.Range('E23').Value = '  +0.28%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.07%  '

$ws.Range('E25').Value = '  -3.83%  '

$ws.Range('E26').Value = '  -0.63%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.10%  '

$ws.Range('E28').Value = '  +3.90%  '

$ws.Range('E29').Value = '  -1.67%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.88%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '163.92'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.17%  '

$ws.Range('E32').Value = '  +2.45%  '

$ws.Range('E33').Value = '  -1.49%  '

$ws.Range('E34').Value = '  +1.14%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.115'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.15%  '

$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.52'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -13.51%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.58'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.54%  '

$ws.Range('E38').Value = '  +0.78%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.68%  '

$ws.Range('E40').Value = '  -6.29%  '

$ws.Range('E41').Value = '  +2.96%  '

$ws.Range('E42').Value = '  +0.20%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '68.33'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.78%  '

$ws.Range('E44').Value = '  -0.93%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -12.79%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.27%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '113.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.18%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '79.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.08%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.29%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.12%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.597.41'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.88%  '
